# Calc.xlsx — add two new skill/version rows (flare_, SpMoonArrow) to the
# "Calc" sheet, matching the latest EA build numbers.
#
# Original data block (rows 3-36):
#   ... row10 ball_ / Alpha 16.1
#       row11 bolt_ / Alpha 16.1
#   ...
#       row15 sword_ / EA 23.117
#       row16 SpMoonSpear / EA 23.191
#   ...
#
# New data block (rows 3-38):
#   ... row10 ball_ / Alpha 16.1
#       row11 flare_ / EA 23.220        <- inserted
#       row12 bolt_ / Alpha 16.1
#   ...
#       row16 sword_ / EA 23.117
#       row17 SpMoonArrow / EA 23.202   <- inserted
#       row18 SpMoonSpear / EA 23.191
#   ...

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert "flare_" right after "ball_" (original row 10 -> new row 11),
# pushing every following row down by one.
$ws.Rows("11").Insert()
$ws.Range("A11").Value = "flare_"
$ws.Range("B11").Value = "EA 23.220"

# Insert "SpMoonArrow" right after "sword_" (now row 16 -> new row 17),
# pushing every following row down by one more.
$ws.Rows("17").Insert()
$ws.Range("A17").Value = "SpMoonArrow"
$ws.Range("B17").Value = "EA 23.202"
